# Applies the "New data Sheets for the Assessment" update:
#  - site_data: fill in newly measured river/catchment values
#  - urban_catchment_landuse: fill in observed land-use proportions and
#    compute the remaining "no_runoff" share with a formula
#  - updates the active sheet/selection/scroll position to match the
#    state the workbook was left in after data entry

$wb = $excel.ActiveWorkbook

# ---- site_data sheet -------------------------------------------------
$siteData = $wb.Worksheets.Item("site_data")

$siteData.Range("D5").Value  = 0.25      # Q_mean [m3/s]
$siteData.Range("D6").Value  = 1         # river_cross_section [m2]
$siteData.Range("D7").Value  = 430       # river_length [m]
$siteData.Range("D13").Value = 56.5      # area_catch [km2]
$siteData.Range("D14").Value = 0.913     # area_urban [km2]
$siteData.Range("D15").Value = 0.913     # area_plan [km2]

# ---- urban_catchment_landuse sheet -----------------------------------
$landuse = $wb.Worksheets.Item("urban_catchment_landuse")

$landuse.Range("C2").Value = 0
$landuse.Range("C3").Value = 0.68
$landuse.Range("C4").Value = 0.07
$landuse.Range("C5").Value = 0
$landuse.Range("C6").Formula = "=1-SUM(C2:C5)"

# ---- restore on-screen selections / active sheet ----------------------
$pollution = $wb.Worksheets.Item("pollution_data")

$landuse.Activate()
$landuse.Range("F20").Select() | Out-Null

$pollution.Activate()
$pollution.Range("D15").Select() | Out-Null

$siteData.Activate()
$siteData.Range("D17").Select() | Out-Null
